$wb = $excel.ActiveWorkbook

# The "Column Relationships" sheet becomes the active/selected tab (it was
# "Column Attributes" before), the cursor/selection on it lands on D18, and
# its data columns (B:R) get widened to a uniform custom width.
$wsRel = $wb.Worksheets.Item("Column Relationships")

$wsRel.Activate()
$wsRel.Range("D18").Select()
$wsRel.Range("B1:R1").ColumnWidth = 9.83
